$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 279 and 280, pushing the existing
# rows (old 279..376) down to become rows 281..378.
$ws.Rows.Item(279).Insert()
$ws.Rows.Item(279).Insert()

# --- New row 279 ---
$ws.Cells.Item(279, 1).Value = 4
$ws.Cells.Item(279, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(279, 3).Value = "Los Lagos"
$ws.Cells.Item(279, 4).Value = 44809
$ws.Cells.Item(279, 5).Value = 10
$ws.Cells.Item(279, 6).Value = "Fruta"
$ws.Cells.Item(279, 7).Value = 100102
$ws.Cells.Item(279, 8).Value = "Cítricos"
$ws.Cells.Item(279, 9).Value = 100102006
$ws.Cells.Item(279, 10).Value = "Pomelo"
$ws.Cells.Item(279, 11).Value = "Start Ruby"
$ws.Cells.Item(279, 12).Value = "Primera"
$ws.Cells.Item(279, 13).Value = 40
$ws.Cells.Item(279, 14).Value = 14000
$ws.Cells.Item(279, 15).Value = 15000
$ws.Cells.Item(279, 16).Value = 14500
$ws.Cells.Item(279, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(279, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(279, 19).Value = 1036
$ws.Cells.Item(279, 20).Value = 14

# --- New row 280 ---
$ws.Cells.Item(280, 1).Value = 4
$ws.Cells.Item(280, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(280, 3).Value = "Los Lagos"
$ws.Cells.Item(280, 4).Value = 44809
$ws.Cells.Item(280, 5).Value = 10
$ws.Cells.Item(280, 6).Value = "Fruta"
$ws.Cells.Item(280, 7).Value = 100102
$ws.Cells.Item(280, 8).Value = "Cítricos"
$ws.Cells.Item(280, 9).Value = 100102006
$ws.Cells.Item(280, 10).Value = "Pomelo"
$ws.Cells.Item(280, 11).Value = "Start Ruby"
$ws.Cells.Item(280, 12).Value = "Segunda"
$ws.Cells.Item(280, 13).Value = 20
$ws.Cells.Item(280, 14).Value = 12000
$ws.Cells.Item(280, 15).Value = 12000
$ws.Cells.Item(280, 16).Value = 12000
$ws.Cells.Item(280, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(280, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(280, 19).Value = 857
$ws.Cells.Item(280, 20).Value = 14
